$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its original text formatting so values
# like "0.9986" or "1.0000" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.579.48'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = '1.672.64'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("D4").Value = '0.9986'
$ws.Range("D5").Value = '239.80'
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4773'
$ws.Range("E7").Value = '  +1.11%  '
$ws.Range("D8").Value = '0.2632'
$ws.Range("E8").Value = '  +2.57%  '
$ws.Range("D9").Value = '0.06183'
$ws.Range("E9").Value = '  +2.82%  '
$ws.Range("D10").Value = '1.673.44'
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("D11").Value = '0.06988'
$ws.Range("E11").Value = '  -1.33%  '
$ws.Range("D12").Value = '14.89'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '0.5918'
$ws.Range("E13").Value = '  -3.97%  '
$ws.Range("D14").Value = '4.386'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '75.40'
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("D16").Value = '1.0000'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = '0.9992'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '25.561.91'
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("D19").Value = '0.000006776'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("D20").Value = '11.44'
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '1.887.02'
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").Value = '4.461'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").Value = '8.771'
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").Value = '5.281'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '136.98'
$ws.Range("E25").Value = '  +3.02%  '
$ws.Range("D26").Value = '15.06'
$ws.Range("E26").Value = '  +1.52%  '
$ws.Range("D27").Value = '1.384'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '1.734'
$ws.Range("E28").Value = '  +4.72%  '
$ws.Range("D29").Value = '104.80'
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("D30").Value = '3.978'
$ws.Range("E30").Value = '  +6.42%  '
$ws.Range("D31").Value = '0.07856'
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("D32").Value = '3.649'
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("D33").Value = '0.9987'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '0.04284'
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").Value = '2.623'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").Value = '0.9587'
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("D37").Value = '0.6088'
$ws.Range("E37").Value = '  +4.53%  '
$ws.Range("D38").Value = '2.588'
$ws.Range("E38").Value = '  +1.43%  '
$ws.Range("D39").Value = '0.8899'
$ws.Range("E39").Value = '  +7.87%  '
$ws.Range("D40").Value = '0.9994'
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = '1.867'
$ws.Range("E41").Value = '  +3.70%  '
$ws.Range("D42").Value = '0.01486'
$ws.Range("E42").Value = '  -4.39%  '
$ws.Range("D43").Value = '96.34'
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("D44").Value = '0.3764'
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("D45").Value = '4.894'
$ws.Range("E45").Value = '  +3.38%  '
$ws.Range("D46").Value = '0.1120'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").Value = '6.241'
$ws.Range("E47").Value = '  +2.48%  '
$ws.Range("D48").Value = '0.05269'
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").Value = '29.96'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").Value = '7.417'
$ws.Range("E50").Value = '  +3.70%  '
$ws.Range("E51").Value = '  +0.29%  '
